# edit.ps1
# Applies the "edit ICD, Database funcs" change:
#   - Rename the "INT"/"int" column-type label to "UINT" across the ICD
#     sheets (user, task, taskHistory, chat, charHistory).
#   - Simplify the taskHistory "command" example rows: the action-type
#     column loses its "...변경/보고" Korean labels in favour of the plain
#     field name, and the command column drops the "[field]" prefix
#     (the due-date example becomes a bare numeric date).
#   - Restore each sheet's remembered selection to match the new state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. ICD type columns: INT / int -> UINT
# ---------------------------------------------------------------------

$wsUser = $wb.Worksheets.Item("user")
$wsUser.Range("D2").Value = "UINT"

$wsTask = $wb.Worksheets.Item("task")
$wsTask.Range("A2").Value = "UINT"
$wsTask.Range("P2").Value = "UINT"
$wsTask.Range("R2").Value = "UINT"

$wsTaskHistory = $wb.Worksheets.Item("taskHistory")
$wsTaskHistory.Range("A2").Value = "UINT"
$wsTaskHistory.Range("B2").Value = "UINT"

$wsChat = $wb.Worksheets.Item("chat")
$wsChat.Range("A2").Value = "UINT"

$wsCharHistory = $wb.Worksheets.Item("charHistory")
$wsCharHistory.Range("A2").Value = "UINT"
$wsCharHistory.Range("B2").Value = "UINT"

# ---------------------------------------------------------------------
# 2. taskHistory example rows (command column) simplified
# ---------------------------------------------------------------------

$wsTaskHistory.Range("C5").Value = "title"
$wsTaskHistory.Range("F5").Value = "B"

$wsTaskHistory.Range("C6").Value = "status"
$wsTaskHistory.Range("F6").Value = "message"

$wsTaskHistory.Range("C7").Value = "priority"
$wsTaskHistory.Range("F7").Value = "urgent"

$wsTaskHistory.Range("C8").Value = "worker"
$wsTaskHistory.Range("F8").Value = "user"

$wsTaskHistory.Range("C9").Value = "due"
$wsTaskHistory.Range("F9").Value = 180506

# ---------------------------------------------------------------------
# 3. Restore per-sheet selections (last Select() wins the active tab,
#    so taskHistory - the originally tab-selected sheet - is done last).
# ---------------------------------------------------------------------

$wsUser.Range("D3").Select()
$wsTask.Range("J3").Select()
$wsChat.Range("A3").Select()
$wsCharHistory.Range("C2").Select()
$wsTaskHistory.Range("F9").Select()
